# 200212 rabbit jump animation modified
# Adds three new to-do rows (25-27) to Sheet1 describing follow-up animation
# work, and updates the sheet's active selection to sit below the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 - new item-usage coding task
$ws.Range("A25").Value = "item사용 결정에 따른 코딩"
$ws.Range("A25").WrapText = $true

# Row 26 - animator note (wraps to 3 lines at the current column width)
$ws.Range("A26").Value = "item눌러서 사용시 애니메이터 추가 필요(현재 sprite변경으로 이루어져서)"
$ws.Range("A26").WrapText = $true
$ws.Rows.Item(26).RowHeight = 49.5

# Row 27 - marble animation / destroy ordering note (wraps to 2 lines)
$ws.Range("A27").Value = "구슬 애니메이션이 destroy되기전에 먼저 이루어지도록 해야함."
$ws.Range("A27").WrapText = $true
$ws.Rows.Item(27).RowHeight = 33

# Move the active selection to just past the newly added rows, matching the
# author's saved cursor position after typing the new entries.
$ws.Range("A28").Select()
